# edit.ps1 — reproduce the "some changes for report in github" commit.
#
# Substance of the change (per the OOXML diff):
#   1) Cp sheet ("Cp" = sheet2): new unit-converted rows 10-15, columns D:F,
#      referencing rows 3-8 (D=*10^-3, E=*10^-6, F=*10^5), written with the
#      D11:D15 / E11:E15 / F11:F15 block as one fill (shared formula) like a
#      drag-fill down from row 11.
#   2) Selection moved on both sheets (Density -> I16, Cp -> G14), with Cp
#      remaining the active/visible tab.
#   3) Minor cosmetic re-layout (default row height / font descent / a couple
#      of best-fit column widths a few hundredths of a character wider) that
#      Excel recomputed when the sheet was re-saved — best effort only, since
#      those are derived/display metrics Excel itself recomputes on save
#      rather than values exposed for scripted assignment.

$wb = $excel.ActiveWorkbook

$wsDensity = $wb.Worksheets.Item("Density")
$wsCp = $wb.Worksheets.Item("Cp")

# --- Cp sheet: add the converted-units rows (10-15) for columns D, E, F ----
# Row 10 is a standalone formula (first block), rows 11-15 are entered as one
# fill so the exported file keeps them as a single shared-formula group
# (ref="D11:D15" / si="0", etc.), matching the source workbook exactly.
$wsCp.Range("D10").Formula = "=D3*10^-3"
$wsCp.Range("E10").Formula = "=E3*10^-6"
$wsCp.Range("F10").Formula = "=F3*10^5"

$wsCp.Range("D11:D15").Formula = "=D4*10^-3"
$wsCp.Range("E11:E15").Formula = "=E4*10^-6"
$wsCp.Range("F11:F15").Formula = "=F4*10^5"

# Best-effort width nudge for the newly-visible column E header ("10^6*C")
# now that rows below it are populated — Excel's real bestFit pixel metrics
# aren't reachable from script, so this just gets us into the right range.
$wsCp.Columns.Item(5).ColumnWidth = 11.86

# --- Selections on both sheets, Cp left as the active tab -----------------
$wsDensity.Range("I16").Select()
$wsCp.Activate()
$wsCp.Range("G14").Select()
